$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the A23 timestamp to a slightly more precise value.
$ws.Range("A23").Value = 44336.77998160533

# Append the newly retrieved data row.
$ws.Range("A24").Value = 44337.78224595184
$ws.Range("B24").Value = 74747
$ws.Range("C24").Value = 62814
$ws.Range("D24").Value = 3161
$ws.Range("E24").Value = 2097
$ws.Range("F24").Value = 1483
$ws.Range("G24").Value = 19417
$ws.Range("H24").Value = 1287
$ws.Range("I24").Value = 845
$ws.Range("J24").Value = 195

# Match the date/time number formatting used for the rest of column A.
$ws.Range("A24").NumberFormat = $ws.Range("A23").NumberFormat
